$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 247, pushing existing rows 247:260 down to 248:261
$ws.Rows("247:247").Insert()

# Populate the newly inserted row 247 with the new weekly data point
$ws.Range("A247").Value = 10
$ws.Range("B247").Value = "Vega Modelo de Temuco"
$ws.Range("C247").Value = "La Araucanía"
$ws.Range("D247").Value = 45021
$ws.Range("E247").Value = 9
$ws.Range("F247").Value = 100112012
$ws.Range("G247").Value = "Espinaca"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 40
$ws.Range("K247").Value = 12000
$ws.Range("L247").Value = 12000
$ws.Range("M247").Value = 12000
$ws.Range("N247").Value = "`$/docena de atados"
$ws.Range("O247").Value = "Región de La Araucanía"
$ws.Range("P247").Value = 4000
$ws.Range("Q247").Value = 3
$ws.Range("R247").Value = "Hortaliza"
